$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain a text value even when the content
    # looks like a number (matches the source inlineStr cells),
    # without leaving a NumberFormat/quotePrefix style artifact behind.
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "26.621.82"
$ws.Range("E2").Value = "  -0.19%  "
Set-TextValue "D3" "1.593.95"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "211.07"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -1.58%  "
Set-TextValue "D10" "19.38"
$ws.Range("E10").Value = "  -1.57%  "
Set-TextValue "D11" "0.0838"
$ws.Range("E11").Value = "  +0.40%  "
Set-TextValue "D12" "1.817.37"
$ws.Range("E12").Value = "  +0.21%  "
Set-TextValue "D13" "1.625.78"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  -1.17%  "
Set-TextValue "D16" "64.49"
$ws.Range("E16").Value = "  -0.10%  "
Set-TextValue "D17" "26.597.13"
$ws.Range("E17").Value = "  -0.19%  "
Set-TextValue "D18" "0.0₃0729"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +0.22%  "
Set-TextValue "D20" "207.44"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("E22").Value = "  -0.11%  "
Set-TextValue "D23" "2.29"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("E24").Value = "  -0.44%  "
Set-TextValue "D25" "145.90"
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  +0.00%  "
Set-TextValue "D27" "7.14"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -0.12%  "
Set-TextValue "D33" "0.655"
$ws.Range("E33").Value = "  +0.47%  "
Set-TextValue "D34" "2.91"
$ws.Range("E34").Value = "  +0.04%  "
Set-TextValue "D35" "1.283.56"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("E40").Value = "  +0.11%  "
Set-TextValue "D41" "5.43"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("E43").Value = "  -0.41%  "
Set-TextValue "D44" "63.46"
$ws.Range("E44").Value = "  -0.03%  "
Set-TextValue "D45" "0.921"
$ws.Range("E45").Value = "  +10.15%  "
Set-TextValue "D46" "1.729.76"
$ws.Range("E46").Value = "  +0.19%  "
Set-TextValue "D47" "89.71"
$ws.Range("E47").Value = "  -0.35%  "
Set-TextValue "D48" "1.59"
$ws.Range("E48").Value = "  -0.38%  "
Set-TextValue "D49" "0.0₆0103"
$ws.Range("E49").Value = "  -2.40%  "
Set-TextValue "D50" "0.101"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("E51").Value = "  -1.46%  "
